$d = $word.ActiveDocument

# The "KEY ACHIEVEMENTS AND IMPACT" section contains an "Impact" sub-heading
# followed by six bullet paragraphs describing job duties. They need to be
# rewritten as four concise, impact-focused accomplishment bullets. We locate
# the bullets by their known text (rather than trusting fixed indices, since
# several of the same phrases also appear earlier in the document under
# "PROFESSIONAL EXPERIENCE"), confirm we found the right ones, then edit them
# in place - replacing/deleting from the bottom of the block upward so that
# earlier paragraph indices stay valid while later ones are removed.

function Find-ParagraphIndex($doc, [string]$needle, [int]$afterIndex) {
    $count = $doc.Paragraphs.Count
    for ($i = $afterIndex + 1; $i -le $count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    throw "Could not find paragraph containing '$needle' after index $afterIndex"
}

function Set-ParagraphText($paragraph, [string]$newText) {
    $rng = $paragraph.Range
    $rng.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude the paragraph mark
    $rng.Text = $newText
}

# Anchor on the "KEY ACHIEVEMENTS AND IMPACT" heading so we operate on the
# correct (second) occurrence of this list, not the similar text earlier in
# the Professional Experience section.
$headingIdx = Find-ParagraphIndex $d "KEY ACHIEVEMENTS AND IMPACT" 0
$impactIdx  = Find-ParagraphIndex $d "Impact" $headingIdx

$idx1 = Find-ParagraphIndex $d "Achieved 87% prediction accuracy" $impactIdx
$idx2 = Find-ParagraphIndex $d "Delivered `$4.9M additional revenue through continuous testing" $idx1
$idx3 = Find-ParagraphIndex $d "Built redistricting platform used by thousands" $idx2
$idx4 = Find-ParagraphIndex $d "Developed longitudinal data analysis methods" $idx3
$idx5 = Find-ParagraphIndex $d "Discovered systematic race coding errors" $idx4
$idx6 = Find-ParagraphIndex $d "Trigonometric algorithm for boundary estimation" $idx5

# Edit bottom-up so earlier indices (idx1..idx3) remain valid.
Set-ParagraphText $d.Paragraphs($idx6) "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
$d.Paragraphs($idx5).Range.Delete()
$d.Paragraphs($idx4).Range.Delete()
Set-ParagraphText $d.Paragraphs($idx3) "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
Set-ParagraphText $d.Paragraphs($idx2) "• 23% conversion rate improvement"
Set-ParagraphText $d.Paragraphs($idx1) "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
